# Auto-update draw results: append the 2025-11-16 Pick 4 draw as a new
# row right after the last existing data row, mirroring the other rows
# (all values stored as literal text, same as the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Force text storage so Excel doesn't auto-coerce the date-looking /
# number-looking values (e.g. "2025-11-16", "251116") into real
# dates/numbers - every other cell in this sheet is stored as text.
$targetRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5))
$targetRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-11-16"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251116"
$ws.Cells.Item($newRow, 4).Value = "8-1-7-0"
$ws.Cells.Item($newRow, 5).Value = "2025-11-16T21:37:03.960+04:00"

# Re-apply the same (default/general) style the rest of the sheet uses,
# so the new row doesn't pick up a stray explicit "Text" number-format
# style - only the cell's stored type is changed, not its look.
$sourceRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 5))
$targetRange.Style = $sourceRange.Style
